# Apply cryptos list update (price/volume refresh + a block reorder in rows 45-49)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for cells whose new value would otherwise be parsed as a number
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D21", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell value updates
$ws.Range("D2").Value = '29.174.91'
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").Value = '1.834.12'
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '240.61'
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("D6").Value = '0.6841'
$ws.Range("E6").Value = '  -1.34%  '
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.3014'
$ws.Range("D9").Value = '0.07478'
$ws.Range("E9").Value = '  -2.66%  '
$ws.Range("D10").Value = '23.09'
$ws.Range("E10").Value = '  -1.11%  '
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("D12").Value = '1.839.57'
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = '5.064'
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("D14").Value = '0.6821'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").Value = '87.00'
$ws.Range("E15").Value = '  -6.70%  '
$ws.Range("D16").Value = '6.187'
$ws.Range("E16").Value = '  -6.06%  '
$ws.Range("D17").Value = '29.182.04'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").Value = '2.080.60'
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("D21").Value = '226.39'
$ws.Range("E21").Value = '  -6.02%  '
$ws.Range("D22").Value = '0.9993'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").Value = '7.425'
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("E25").Value = '  -3.23%  '
$ws.Range("D26").Value = '159.92'
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").Value = '8.761'
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").Value = '18.09'
$ws.Range("E28").Value = '  -0.40%  '
$ws.Range("D29").Value = '1.507'
$ws.Range("E29").Value = '  -2.19%  '
$ws.Range("D30").Value = '4.257'
$ws.Range("E30").Value = '  +1.09%  '
$ws.Range("D31").Value = '4.143'
$ws.Range("E31").Value = '  -0.86%  '
$ws.Range("D32").Value = '1.202'
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("D33").Value = '0.05160'
$ws.Range("E33").Value = '  +1.17%  '
$ws.Range("D34").Value = '0.7666'
$ws.Range("E34").Value = '  -1.83%  '
$ws.Range("D35").Value = '1.839'
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("D37").Value = '2.671'
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("D38").Value = '1.306.55'
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("D39").Value = '0.01833'
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("D40").Value = '2.722'
$ws.Range("E40").Value = '  +0.66%  '
$ws.Range("D41").Value = '0.9364'
$ws.Range("E41").Value = '  -1.72%  '
$ws.Range("D42").Value = '5.809'
$ws.Range("E42").Value = '  -5.49%  '
$ws.Range("D43").Value = '104.35'
$ws.Range("E43").Value = '  -2.47%  '
$ws.Range("D44").Value = '0.9992'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '9.612'
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.982.19'
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.5205'
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.00000000123'
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '64.89'
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("D50").Value = '1.772'
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").Value = '0.05917'
$ws.Range("E51").Value = '  +1.00%  '
